$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing data before restructuring -----------------------
# Old layout (row1 headers / row2 values), columns A-J:
#   A-H : LowWait, MediumWait, HighWait, HighestWait,
#         LowExplicitWait, MediumExplicitWait, LongExplicitWait, HighestExplicitWait
#   I   : CriteriaSelect / "GNB AutoPost Criiteria Set"
#   J   : SubmissionNotes / "Testing2920204"
# Columns K (blank), L (URL), M (UserName), N (Password), O, P (blank) are dropped.

$waitHeaders = @()
$waitValues = @()
for ($col = 1; $col -le 8; $col++) {
    $waitHeaders += $ws.Cells.Item(1, $col).Value2
    $waitValues += $ws.Cells.Item(2, $col).Value2
}
$criteriaHeader = $ws.Cells.Item(1, 9).Value2
$criteriaValue = $ws.Cells.Item(2, 9).Value2
$notesHeader = $ws.Cells.Item(1, 10).Value2
$notesValue = $ws.Cells.Item(2, 10).Value2

# --- Wipe the sheet and rebuild with the new column order ------------------
$ws.Cells.Clear()

# New column order: CriteriaSelect, SubmissionNotes, then the eight wait columns, then a blank column.
$newHeaders = @($criteriaHeader, $notesHeader) + $waitHeaders
$newValues = @($criteriaValue, $notesValue) + $waitValues

for ($col = 1; $col -le $newHeaders.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $newHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.Borders.LineStyle = 1
}
# Trailing blank column (K) keeps the header border/bold formatting, no value.
$lastCol = $newHeaders.Length + 1
$blankHeaderCell = $ws.Cells.Item(1, $lastCol)
$blankHeaderCell.Font.Bold = $true
$blankHeaderCell.HorizontalAlignment = -4108
$blankHeaderCell.Borders.LineStyle = 1

for ($col = 1; $col -le $newValues.Length; $col++) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.Value = $newValues[$col - 1]
    $cell.Borders.LineStyle = 1
}
# Trailing blank column (K) on the data row, border only.
$ws.Cells.Item(2, $lastCol).Borders.LineStyle = 1

# --- Column widths (best-fit to content) -----------------------------------
$ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(2, $lastCol)).Columns.AutoFit()

# --- Remove the now-unused "Hyperlink" cell style ---------------------------
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "Hyperlink") {
        $s.Delete()
    }
}

# --- Sheet view / selection -------------------------------------------------
$ws.Range("D9").Select()

# --- Workbook-level metadata -------------------------------------------------
$wb.Application.AlertBeforeOverwriting = $false
